$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 305, pushing the existing rows 305-308
# (and everything below) down to 309-312.
$ws.Rows.Item(305).Insert()
$ws.Rows.Item(305).Insert()
$ws.Rows.Item(305).Insert()
$ws.Rows.Item(305).Insert()

# Populate the 4 newly inserted rows (305-308) with the new weekly
# Melón price entries for Feria Lagunitas de Puerto Montt.
$rows = @(
    @{ Row = 305; D = 44610; H = "Calameño"; I = "Extra";   J = 3000; K = 1500; L = 1500; M = 1500; N = "`$/unidad"; P = 1500; Q = 1 },
    @{ Row = 306; D = 44610; H = "Calameño"; I = "Primera"; J = 3000; K = 1200; L = 1200; M = 1200; N = "`$/unidad"; P = 1200; Q = 1 },
    @{ Row = 307; D = 44610; H = "Tuna";     I = "Extra";   J = 3000; K = 1500; L = 1500; M = 1500; N = "`$/unidad"; P = 1500; Q = 1 },
    @{ Row = 308; D = 44610; H = "Tuna";     I = "Primera"; J = 3000; K = 1200; L = 1200; M = 1200; N = "`$/unidad"; P = 1200; Q = 1 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 4
    $ws.Cells.Item($row, 2).Value = "Feria Lagunitas de Puerto Montt"
    $ws.Cells.Item($row, 3).Value = "Los Lagos"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 10
    $ws.Cells.Item($row, 6).Value = 100112027
    $ws.Cells.Item($row, 7).Value = "Melón"
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = "Región de O'Higgins"
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
